$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C ("UL_DESC"), shifting MFR/MFR_PARTNUM/etc. right.
$ws.Columns.Item(3).Insert()

# Header
$ws.Range("C1").Value = "UL_DESC"

# New UL_DESC values for each data row (matches the UL part number / voltage pattern)
$ws.Range("C2").Value = "VARIABLE FREQUENCY DRIVE, 230V, 5 A"
$ws.Range("C3").Value = "VARIABLE FREQUENCY DRIVE, 460V, 3.4 A"
$ws.Range("C4").Value = "VARIABLE FREQUENCY DRIVE, 230V, 8 A"
$ws.Range("C5").Value = "VARIABLE FREQUENCY DRIVE, 460V, 4.8 A"
$ws.Range("C6").Value = "VARIABLE FREQUENCY DRIVE, 230V, 17.5 A"
$ws.Range("C7").Value = "VARIABLE FREQUENCY DRIVE, 460V, 9.2 A"
$ws.Range("C8").Value = "VARIABLE FREQUENCY DRIVE, 230V, 25 A"
$ws.Range("C9").Value = "VARIABLE FREQUENCY DRIVE, 460V, 14.8 A"

# Resize the new column to fit its (longer) contents
$ws.Columns.Item(3).ColumnWidth = 38

# Move the selection, as recorded in the saved workbook view
$ws.Range("C15").Select()

# Reflect the window's new horizontal screen position, as recorded in the workbook view
$excel.Windows.Item(1).Left = 3780
